$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 653.2222
$ws.Range("I5").Value = 13.166667
$ws.Range("J5").Value = 1933.3334
$ws.Range("K5").Value = 13.166667
$ws.Range("L5").Value = 1933.3334
$ws.Range("M5").Value = 101.833333
$ws.Range("N5").Value = -2163.3334
$ws.Range("H15").Value = 218.38
$ws.Range("I15").Value = 218.38
$ws.Range("K15").Value = 655.14
$ws.Range("M15").Value = -486.14
$ws.Range("H53").Value = 3857
$ws.Range("I53").Value = 392.5
$ws.Range("J53").Value = 5836.7144
$ws.Range("K53").Value = 392.5
$ws.Range("L53").Value = 5836.7144
$ws.Range("M53").Value = 244.5
$ws.Range("N53").Value = -7110.7144
$ws.Range("H113").Value = 71432296
$ws.Range("I113").Value = 250001280
$ws.Range("J113").Value = 4699.5
$ws.Range("K113").Value = 250001280
$ws.Range("L113").Value = 4699.5
$ws.Range("M113").Value = -249998026
$ws.Range("N113").Value = -11207.5
$ws.Range("H129").Value = 124365.53
$ws.Range("J129").Value = 141847.03
$ws.Range("L129").Value = 425541.09
$ws.Range("N129").Value = -435541.09
$ws.Range("H132").Value = 3793.923
$ws.Range("I132").Value = 4131.522
$ws.Range("K132").Value = 12394.566
$ws.Range("M132").Value = -9864.565999999999
$ws.Range("H137").Value = 1548.3636
$ws.Range("I137").Value = 1249.125
$ws.Range("K137").Value = 3747.375
$ws.Range("M137").Value = -1197.375
$ws.Range("H138").Value = 2450.2144
$ws.Range("I138").Value = 2041.6666
$ws.Range("J138").Value = 2534.7415
$ws.Range("K138").Value = 6124.9998
$ws.Range("L138").Value = 7604.2245
$ws.Range("M138").Value = -984.9997999999996
$ws.Range("N138").Value = -17884.2245
$ws.Range("H139").Value = 50620
$ws.Range("J139").Value = 50620
$ws.Range("L139").Value = 50620
$ws.Range("N139").Value = -60900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 828.97144
$ws.Range("I2").Value = 760.6
$ws.Range("J2").Value = 999.9
$ws.Range("K2").Value = 760.6
$ws.Range("L2").Value = 999.9
$ws.Range("M2").Value = -647.6
$ws.Range("N2").Value = -1225.9
$ws.Range("H32").Value = 5880.0386
$ws.Range("I32").Value = 4491.8696
$ws.Range("K32").Value = 4491.8696
$ws.Range("M32").Value = -4204.8696
$ws.Range("H74").Value = 18519646
$ws.Range("I74").Value = 22222712
$ws.Range("K74").Value = 22222712
$ws.Range("M74").Value = -22221838
$ws.Range("H77").Value = 18519646
$ws.Range("I77").Value = 22222712
$ws.Range("K77").Value = 111113560
$ws.Range("M77").Value = -111109192
$ws.Range("H116").Value = 828.97144
$ws.Range("I116").Value = 760.6
$ws.Range("J116").Value = 999.9
$ws.Range("K116").Value = 760.6
$ws.Range("L116").Value = 999.9
$ws.Range("M116").Value = 1533.4
$ws.Range("N116").Value = -5587.9
$ws.Range("H122").Value = 1564.84
$ws.Range("I122").Value = 1533.6
$ws.Range("K122").Value = 4600.799999999999
$ws.Range("M122").Value = -2150.799999999999
$ws.Range("H132").Value = 14373.585
$ws.Range("I132").Value = 1938.4193
$ws.Range("J132").Value = 52922.6
$ws.Range("K132").Value = 5815.257900000001
$ws.Range("L132").Value = 158767.8
$ws.Range("M132").Value = -3285.257900000001
$ws.Range("N132").Value = -163827.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 828.97144
$ws.Range("I3").Value = 760.6
$ws.Range("J3").Value = 999.9
$ws.Range("K3").Value = 760.6
$ws.Range("L3").Value = 999.9
$ws.Range("M3").Value = -646.6
$ws.Range("N3").Value = -1227.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3509.0715
$ws.Range("I31").Value = 1819.6666
$ws.Range("K31").Value = 1819.6666
$ws.Range("M31").Value = -1524.6666
$ws.Range("H34").Value = 3509.0715
$ws.Range("I34").Value = 1819.6666
$ws.Range("K34").Value = 1819.6666
$ws.Range("M34").Value = -1617.6666
$ws.Range("H132").Value = 3056.5
$ws.Range("I132").Value = 2180.9546
$ws.Range("J132").Value = 6266.8335
$ws.Range("K132").Value = 6542.8638
$ws.Range("L132").Value = 18800.5005
$ws.Range("M132").Value = -4012.8638
$ws.Range("N132").Value = -23860.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1292.6364
$ws.Range("I5").Value = 904.37933
$ws.Range("K5").Value = 2713.13799
$ws.Range("M5").Value = -2601.13799
$ws.Range("H13").Value = 375
$ws.Range("I13").Value = 150
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = 450
$ws.Range("L13").Value = 1800
$ws.Range("M13").Value = -282
$ws.Range("N13").Value = -2136
$ws.Range("H59").Value = 2816.6667
$ws.Range("I59").Value = 800
$ws.Range("K59").Value = 2400
$ws.Range("M59").Value = -1860
$ws.Range("H81").Value = 5409.9165
$ws.Range("J81").Value = 5409.9165
$ws.Range("L81").Value = 16229.7495
$ws.Range("N81").Value = -18475.7495
$ws.Range("H84").Value = 5409.9165
$ws.Range("J84").Value = 5409.9165
$ws.Range("L84").Value = 48689.2485
$ws.Range("N84").Value = -59921.2485
$ws.Range("H131").Value = 716.12244
$ws.Range("I131").Value = 428.33334
$ws.Range("J131").Value = 734.8913
$ws.Range("K131").Value = 1285.00002
$ws.Range("L131").Value = 2204.6739
$ws.Range("M131").Value = 3754.99998
$ws.Range("N131").Value = -12284.6739
$ws.Range("H135").Value = 1292.6364
$ws.Range("I135").Value = 904.37933
$ws.Range("K135").Value = 8139.41397
$ws.Range("M135").Value = -5604.41397
$ws.Range("H136").Value = 3535.7856
$ws.Range("J136").Value = 4996.778
$ws.Range("L136").Value = 14990.334
$ws.Range("N136").Value = -25190.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2598.8684
$ws.Range("I122").Value = 1951.9231
$ws.Range("J122").Value = 4000.5833
$ws.Range("K122").Value = 5855.7693
$ws.Range("L122").Value = 12001.7499
$ws.Range("M122").Value = -3405.7693
$ws.Range("N122").Value = -16901.7499
$ws.Range("H126").Value = 3919.0637
$ws.Range("I126").Value = 4408.048
$ws.Range("J126").Value = 3524.1155
$ws.Range("K126").Value = 13224.144
$ws.Range("L126").Value = 10572.3465
$ws.Range("M126").Value = -10754.144
$ws.Range("N126").Value = -15512.3465
$ws.Range("H132").Value = 21218.555
$ws.Range("I132").Value = 2177.7222
$ws.Range("J132").Value = 59300.223
$ws.Range("K132").Value = 6533.1666
$ws.Range("L132").Value = 177900.669
$ws.Range("M132").Value = -4003.1666
$ws.Range("N132").Value = -182960.669
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3311.68
$ws.Range("I7").Value = 2955.4375
$ws.Range("J7").Value = 3945
$ws.Range("K7").Value = 2955.4375
$ws.Range("L7").Value = 3945
$ws.Range("M7").Value = -2843.4375
$ws.Range("N7").Value = -4169
$ws.Range("H22").Value = 2620.739
$ws.Range("J22").Value = 1126.4445
$ws.Range("L22").Value = 1126.4445
$ws.Range("N22").Value = -1716.4445
$ws.Range("H27").Value = 2620.739
$ws.Range("J27").Value = 1126.4445
$ws.Range("L27").Value = 1126.4445
$ws.Range("N27").Value = -1340.4445
$ws.Range("H40").Value = 2924.7188
$ws.Range("I40").Value = 2700.04
$ws.Range("K40").Value = 2700.04
$ws.Range("M40").Value = -2564.04
$ws.Range("H46").Value = 1659.5938
$ws.Range("I46").Value = 1817.9445
$ws.Range("K46").Value = 1817.9445
$ws.Range("M46").Value = -1629.9445
$ws.Range("H55").Value = 1255.25
$ws.Range("I55").Value = 1591.6666
$ws.Range("K55").Value = 1591.6666
$ws.Range("M55").Value = -1418.6666
$ws.Range("H122").Value = 579823.4399999999
$ws.Range("I122").Value = 894062.5600000001
$ws.Range("K122").Value = 2682187.68
$ws.Range("M122").Value = -2679737.68
$ws.Range("H126").Value = 3311.68
$ws.Range("I126").Value = 2955.4375
$ws.Range("J126").Value = 3945
$ws.Range("K126").Value = 8866.3125
$ws.Range("L126").Value = 11835
$ws.Range("M126").Value = -6396.3125
$ws.Range("N126").Value = -16775
$ws.Range("H132").Value = 485075.9
$ws.Range("I132").Value = 862506.8
$ws.Range("J132").Value = 4709.364
$ws.Range("K132").Value = 2587520.4
$ws.Range("L132").Value = 14128.092
$ws.Range("M132").Value = -2584990.4
$ws.Range("N132").Value = -19188.092
$ws.Range("H136").Value = 1811.3043
$ws.Range("I136").Value = 1602.8572
$ws.Range("K136").Value = 4808.571599999999
$ws.Range("M136").Value = -2258.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4003
$ws.Range("J62").Value = 4003
$ws.Range("L62").Value = 4003
$ws.Range("N62").Value = -5251
$ws.Range("H65").Value = 4003
$ws.Range("J65").Value = 4003
$ws.Range("L65").Value = 20015
$ws.Range("N65").Value = -26255
$ws.Range("H101").Value = 16175
$ws.Range("J101").Value = 16175
$ws.Range("L101").Value = 16175
$ws.Range("N101").Value = -22665
$ws.Range("H126").Value = 1669.3334
$ws.Range("I126").Value = 1257.5
$ws.Range("J126").Value = 3316.6667
$ws.Range("K126").Value = 3772.5
$ws.Range("L126").Value = 9950.000100000001
$ws.Range("M126").Value = -1302.5
$ws.Range("N126").Value = -14890.0001
$ws.Range("H132").Value = 1993.36
$ws.Range("I132").Value = 1372.7059
$ws.Range("K132").Value = 4118.1177
$ws.Range("M132").Value = -1588.1177
$ws.Range("H136").Value = 32261822
$ws.Range("I136").Value = 46922596
$ws.Range("J136").Value = 8120.5
$ws.Range("K136").Value = 140767788
$ws.Range("L136").Value = 24361.5
$ws.Range("M136").Value = -140765238
